$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "nome" in C1, matching the bold/centered/bordered
#     style already used by A1/B1 (style index carried via copy/paste of
#     formats only, so no new style gets created). ---
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1").Value = "nome"

# --- New row 2: numero / codigo / nome for "andy" ---
# The numero/codigo columns hold digit-only strings that must stay text
# (not be auto-converted to numbers). Temporarily force a text format,
# assign the value, then drop the format again so the cell ends up using
# the plain default style (same as the untouched data cells).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "48998418335"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").ClearFormats()

$ws.Range("C2").Value = "andy "

# --- Row 3: fix codigo (was "1.0", now "1") and add nome "sla" ---
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1"
$ws.Range("B3").ClearFormats()

$ws.Range("C3").Value = "sla"
